$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    "S-2021566152",
    "S-202191694",
    "S-2021578905",
    "S-2021905948",
    "S-2021464411",
    "S-2021632681",
    "S-2021897308",
    "S-2021462804",
    "S-2021215360",
    "S-202181321",
    "S-2021224053",
    "S-2021585403",
    "S-2021838439",
    "S-2021613471",
    "S-2021452415",
    "S-2021471410"
)

$startRow = 29
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
